$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
    $tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U79"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
    Write-Host "table added: " $tbl.Name
} catch {
    Write-Host "ERROR: $_"
}
